# Update latest output (run 126)
# Applies updated optimisation results to the "Schedule" and "Detailed" sheets

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (Cost ($) and Unit Cost ($/ML) columns) ---
$schedule.Range("E3").Value = -2.811704999999982
$schedule.Range("F3").Value = -0.0826485890652552
$schedule.Range("E4").Value = 575.5282845
$schedule.Range("F4").Value = 30.45123198412699
$schedule.Range("E5").Value = -185.71424625
$schedule.Range("F5").Value = -5.458972552910054

# --- Detailed sheet updates (Price and Type columns) ---
$detailed.Range("B30").Value = -9.43817
$detailed.Range("B31").Value = -9.5
$detailed.Range("B32").Value = 14.00051
$detailed.Range("C32").Value = "historical"
$detailed.Range("B33").Value = 0.00048
$detailed.Range("C33").Value = "historical"
$detailed.Range("B34").Value = 10.31935
$detailed.Range("B35").Value = 22.07
$detailed.Range("B36").Value = 25.73712
$detailed.Range("B37").Value = 49.80515
$detailed.Range("B38").Value = 57.3
$detailed.Range("B39").Value = 63.92984
$detailed.Range("B40").Value = 65
$detailed.Range("B41").Value = 66.16136
$detailed.Range("B42").Value = 73.20007
$detailed.Range("B43").Value = 71.92849
$detailed.Range("B44").Value = 68.79257
$detailed.Range("B45").Value = 63.47425
$detailed.Range("B47").Value = 63.96083
$detailed.Range("B48").Value = 62.40305
$detailed.Range("B49").Value = 61.34585
$detailed.Range("B50").Value = 59.97554
$detailed.Range("B59").Value = 67.70496
$detailed.Range("B61").Value = 79.95022
$detailed.Range("B62").Value = 79.95016
$detailed.Range("B63").Value = 67.70788
$detailed.Range("B66").Value = 0.51002
$detailed.Range("B67").Value = -2.83936
$detailed.Range("B68").Value = -6.54017
$detailed.Range("B69").Value = -6.96691
$detailed.Range("B70").Value = -8.566689999999999
$detailed.Range("B71").Value = -9.75165
$detailed.Range("B72").Value = -15.60246
$detailed.Range("B73").Value = -22.35626
$detailed.Range("B74").Value = -22.90284
$detailed.Range("B75").Value = -23.5
$detailed.Range("B76").Value = -23.5
$detailed.Range("B77").Value = -27
$detailed.Range("B78").Value = -23.93719
$detailed.Range("B79").Value = -23.13936
$detailed.Range("B80").Value = -22.06328
$detailed.Range("B81").Value = -10
$detailed.Range("B82").Value = -5.74405
$detailed.Range("B83").Value = -2.03998
$detailed.Range("B85").Value = 46.3101
$detailed.Range("B88").Value = 78
$detailed.Range("B89").Value = 105.0001
$detailed.Range("B90").Value = 86.97678000000001
$detailed.Range("B91").Value = 73.75449
$detailed.Range("B92").Value = 70.0634
$detailed.Range("B93").Value = 73.19
$detailed.Range("B94").Value = 60.2421
$detailed.Range("B95").Value = 65
$detailed.Range("B96").Value = 64.8901
$detailed.Range("B97").Value = 64.8901
